$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 535
$ws.Range("I61").Value = 535
$ws.Range("K61").Value = 1605
$ws.Range("M61").Value = -1433
$ws.Range("H70").Value = 3187.074
$ws.Range("I70").Value = 3122.4
$ws.Range("K70").Value = 9367.200000000001
$ws.Range("M70").Value = -9097.200000000001
$ws.Range("H73").Value = 3187.074
$ws.Range("I73").Value = 3122.4
$ws.Range("K73").Value = 9367.200000000001
$ws.Range("M73").Value = -8431.200000000001
$ws.Range("H87").Value = 19999.87
$ws.Range("J87").Value = 19999.87
$ws.Range("L87").Value = 19999.87
$ws.Range("N87").Value = -22495.87
$ws.Range("H90").Value = 19999.87
$ws.Range("J90").Value = 19999.87
$ws.Range("L90").Value = 59999.61
$ws.Range("N90").Value = -72479.61
$ws.Range("H113").Value = 4754.1763
$ws.Range("I113").Value = 3930.5
$ws.Range("J113").Value = 5486.3335
$ws.Range("K113").Value = 3930.5
$ws.Range("L113").Value = 5486.3335
$ws.Range("M113").Value = -676.5
$ws.Range("N113").Value = -11994.3335
$ws.Range("H125").Value = 1399.3334
$ws.Range("I125").Value = 1293.2
$ws.Range("J125").Value = 1532
$ws.Range("K125").Value = 11638.8
$ws.Range("L125").Value = 13788
$ws.Range("M125").Value = -9178.800000000001
$ws.Range("N125").Value = -18708

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10081.25
$ws.Range("I32").Value = 9123
$ws.Range("K32").Value = 9123
$ws.Range("M32").Value = -8836
$ws.Range("H61").Value = 3164.913
$ws.Range("I61").Value = 1932.6666
$ws.Range("J61").Value = 3957.0715
$ws.Range("K61").Value = 1932.6666
$ws.Range("L61").Value = 3957.0715
$ws.Range("M61").Value = -1720.6666
$ws.Range("N61").Value = -4381.0715
$ws.Range("H74").Value = 1562.9048
$ws.Range("I74").Value = 1511.6842
$ws.Range("K74").Value = 1511.6842
$ws.Range("M74").Value = -637.6841999999999
$ws.Range("H77").Value = 1562.9048
$ws.Range("I77").Value = 1511.6842
$ws.Range("K77").Value = 7558.420999999999
$ws.Range("M77").Value = -3190.420999999999
$ws.Range("H110").Value = 3585.4443
$ws.Range("I110").Value = 3502.2354
$ws.Range("K110").Value = 3502.2354
$ws.Range("M110").Value = -1457.2354
$ws.Range("H122").Value = 3832.7437
$ws.Range("I122").Value = 1848.3334
$ws.Range("K122").Value = 5545.0002
$ws.Range("M122").Value = -3095.0002
$ws.Range("H125").Value = 68955.625
$ws.Range("J125").Value = 68955.625
$ws.Range("L125").Value = 68955.625
$ws.Range("N125").Value = -78795.625
$ws.Range("H132").Value = 3076.0386
$ws.Range("I132").Value = 2833.6584
$ws.Range("K132").Value = 8500.975199999999
$ws.Range("M132").Value = -5970.975199999999
$ws.Range("H136").Value = 3164.913
$ws.Range("I136").Value = 1932.6666
$ws.Range("J136").Value = 3957.0715
$ws.Range("K136").Value = 5797.9998
$ws.Range("L136").Value = 11871.2145
$ws.Range("M136").Value = -3247.9998
$ws.Range("N136").Value = -16971.2145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 251
$ws.Range("J7").Value = 404.83334
$ws.Range("L7").Value = 404.83334
$ws.Range("N7").Value = -630.83334
$ws.Range("H31").Value = 2885.2173
$ws.Range("I31").Value = 1418.2354
$ws.Range("J31").Value = 7041.6665
$ws.Range("K31").Value = 1418.2354
$ws.Range("L31").Value = 7041.6665
$ws.Range("M31").Value = -1123.2354
$ws.Range("N31").Value = -7631.6665
$ws.Range("H34").Value = 2885.2173
$ws.Range("I34").Value = 1418.2354
$ws.Range("J34").Value = 7041.6665
$ws.Range("K34").Value = 1418.2354
$ws.Range("L34").Value = 7041.6665
$ws.Range("M34").Value = -1216.2354
$ws.Range("N34").Value = -7445.6665
$ws.Range("H58").Value = 1833.7778
$ws.Range("I58").Value = 1323.9231
$ws.Range("J58").Value = 3159.4
$ws.Range("K58").Value = 1323.9231
$ws.Range("L58").Value = 3159.4
$ws.Range("M58").Value = -1120.9231
$ws.Range("N58").Value = -3565.4
$ws.Range("H122").Value = 1281346.4
$ws.Range("I122").Value = 3403923.8
$ws.Range("K122").Value = 10211771.4
$ws.Range("M122").Value = -10209321.4
$ws.Range("H132").Value = 2983.9062
$ws.Range("I132").Value = 2326.1035
$ws.Range("K132").Value = 6978.310500000001
$ws.Range("M132").Value = -4448.310500000001
$ws.Range("H136").Value = 1833.7778
$ws.Range("I136").Value = 1323.9231
$ws.Range("J136").Value = 3159.4
$ws.Range("K136").Value = 3971.7693
$ws.Range("L136").Value = 9478.200000000001
$ws.Range("M136").Value = -1421.7693
$ws.Range("N136").Value = -14578.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 157.5
$ws.Range("J12").Value = 173.60869
$ws.Range("L12").Value = 520.82607
$ws.Range("N12").Value = -866.82607
$ws.Range("H80").Value = 5832.5
$ws.Range("I80").Value = 5248.75
$ws.Range("K80").Value = 15746.25
$ws.Range("M80").Value = -14810.25
$ws.Range("H83").Value = 5832.5
$ws.Range("I83").Value = 5248.75
$ws.Range("K83").Value = 47238.75
$ws.Range("M83").Value = -42558.75
$ws.Range("H92").Value = 640
$ws.Range("I92").Value = 463.5
$ws.Range("J92").Value = 728.25
$ws.Range("K92").Value = 1390.5
$ws.Range("L92").Value = 2184.75
$ws.Range("M92").Value = -142.5
$ws.Range("N92").Value = -4680.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3955.0952
$ws.Range("I126").Value = 2329.6667
$ws.Range("K126").Value = 6989.000100000001
$ws.Range("M126").Value = -4519.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2696.75
$ws.Range("I7").Value = 1454.7142
$ws.Range("K7").Value = 1454.7142
$ws.Range("M7").Value = -1342.7142
$ws.Range("H22").Value = 1114.125
$ws.Range("I22").Value = 916.1429000000001
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 916.1429000000001
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -621.1429000000001
$ws.Range("N22").Value = -3090
$ws.Range("H27").Value = 1114.125
$ws.Range("I27").Value = 916.1429000000001
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 916.1429000000001
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -809.1429000000001
$ws.Range("N27").Value = -2714
$ws.Range("H40").Value = 7801.8438
$ws.Range("I40").Value = 8290.941000000001
$ws.Range("K40").Value = 8290.941000000001
$ws.Range("M40").Value = -8154.941000000001
$ws.Range("H46").Value = 8449.888999999999
$ws.Range("I46").Value = 2171
$ws.Range("J46").Value = 9965.482
$ws.Range("K46").Value = 2171
$ws.Range("L46").Value = 9965.482
$ws.Range("M46").Value = -1983
$ws.Range("N46").Value = -10341.482
$ws.Range("H61").Value = 2096.1562
$ws.Range("I61").Value = 733.7692
$ws.Range("J61").Value = 7999.8335
$ws.Range("K61").Value = 733.7692
$ws.Range("L61").Value = 7999.8335
$ws.Range("M61").Value = -531.7692
$ws.Range("N61").Value = -8403.833500000001
$ws.Range("H100").Value = 83700.21000000001
$ws.Range("I100").Value = 223900.8
$ws.Range("K100").Value = 223900.8
$ws.Range("M100").Value = -223359.8
$ws.Range("H113").Value = 2096.1562
$ws.Range("I113").Value = 733.7692
$ws.Range("J113").Value = 7999.8335
$ws.Range("K113").Value = 733.7692
$ws.Range("L113").Value = 7999.8335
$ws.Range("M113").Value = 1436.2308
$ws.Range("N113").Value = -12339.8335
$ws.Range("H126").Value = 2696.75
$ws.Range("I126").Value = 1454.7142
$ws.Range("K126").Value = 4364.142599999999
$ws.Range("M126").Value = -1894.142599999999
$ws.Range("H132").Value = 3496.3914
$ws.Range("I132").Value = 2862.4722
$ws.Range("K132").Value = 8587.4166
$ws.Range("M132").Value = -6057.4166

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 406.52942
$ws.Range("I107").Value = 387.45456
$ws.Range("J107").Value = 441.5
$ws.Range("K107").Value = 1162.36368
$ws.Range("L107").Value = 1324.5
$ws.Range("M107").Value = 757.6363200000001
$ws.Range("N107").Value = -5164.5
$ws.Range("H132").Value = 2700.1914
$ws.Range("I132").Value = 2526.024
$ws.Range("J132").Value = 4163.2
$ws.Range("K132").Value = 7578.072
$ws.Range("L132").Value = 12489.6
$ws.Range("M132").Value = -5048.072
$ws.Range("N132").Value = -17549.6
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140
$ws.Range("H136").Value = 2580.2942
$ws.Range("I136").Value = 1279.24
$ws.Range("K136").Value = 3837.72
$ws.Range("M136").Value = -1287.72

Write-Output "Applied 216 cell updates across 7 sheets"